{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Goal (per the commit diff):\n//  - Rewrite the \"Jeliko\u017e se jedn\u00e1 o anal\u00fdzu...\" paragraph (previously split\n//    across 4 runs) into a single run with revised/expanded Czech wording.\n//  - Add a new paragraph right after it with an additional concluding\n//    sentence.\n//  - Add a new, empty paragraph right after that (before the section break).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the target paragraph robustly (by its distinctive original text)\n// instead of assuming a fixed index.\nconst marker = \"Jeliko\u017e se jedn\u00e1\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the target paragraph (text starting with '\" + marker + \"').\");\n}\n\nconst newFirstParagraphText =\n  \"Vzhledem k tomu, \u017ee se jedn\u00e1 o anal\u00fdzu za v\u00edce let, byly pou\u017eity pr\u016fm\u011brn\u00e9 ro\u010dn\u00ed hodnoty, aby se zajistila stabilita a p\u0159esnost v\u00fdsledk\u016f. Uva\u017eovali jsme tak\u00e9 o mo\u017enosti \u010derpat ze \u010dtvrtletn\u00edch dat, av\u0161ak tato varianta byla odm\u00edtnuta kv\u016fli sez\u00f3nn\u00edm v\u00fdkyv\u016fm v cen\u00e1ch potravin. Pou\u017eit\u00ed ro\u010dn\u00edch pr\u016fm\u011br\u016f t\u00edmto zp\u016fsobem eliminovalo sez\u00f3nn\u00ed vlivy a poskytlo \u010dist\u0161\u00ed data pro anal\u00fdzu dlouhodob\u00fdch trend\u016f.\";\n\nconst newSecondParagraphText =\n  \"Tento p\u0159\u00edstup n\u00e1m umo\u017enil l\u00e9pe porozum\u011bt dynamice v\u00fdvoje mezd a cen potravin v kontextu ekonomick\u00e9ho prost\u0159ed\u00ed \u010cesk\u00e9 republiky.\";\n\n// Replace the paragraph's (multi-run) content with a single run of the\n// revised text.\ntarget.clear();\ntarget.insertText(newFirstParagraphText, Word.InsertLocation.start);\n\n// Insert the new concluding paragraph right after it.\nconst secondParagraph = target.insertParagraph(newSecondParagraphText, Word.InsertLocation.after);\n\n// Insert a trailing empty paragraph right after that one.\nsecondParagraph.insertParagraph(\"\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Goal (per the commit diff):\n#  - Rewrite the \"Jeliko\u017e se jedn\u00e1 o anal\u00fdzu...\" paragraph (previously split\n#    across 4 runs) into a single run with revised/expanded Czech wording.\n#  - Add a new paragraph right after it with an additional concluding\n#    sentence.\n#  - Add a new, empty paragraph right after that (before the section break).\n\n$d = $word.ActiveDocument\n\n$marker = \"Jeliko\u017e se jedn\u00e1\"\n$oldText = \"Jeliko\u017e se jedn\u00e1 o anal\u00fdzu za v\u00edce let, byly pou\u017eity pr\u016fm\u011brn\u00e9 ro\u010dn\u00ed hodnoty. V \u00favahu byla vzata i mo\u017enost \u010derpat ze \u010dtvrtletn\u00edch dat, ale u cen potravin by to nebylo vhodn\u00e9 kv\u016fli sez\u00f3nn\u00edm v\u00fdkyv\u016fm. T\u00edmto krokem byly nav\u00edc eliminov\u00e1ny sez\u00f3nn\u00ed vlivy. P\u0159i zpracov\u00e1n\u00ed dat byl kladen d\u016fraz na to, aby bylo co nejv\u00edce viditeln\u00fdch informac\u00ed na jedn\u00e9 str\u00e1nce.\"\n$newFirstParagraphText = \"Vzhledem k tomu, \u017ee se jedn\u00e1 o anal\u00fdzu za v\u00edce let, byly pou\u017eity pr\u016fm\u011brn\u00e9 ro\u010dn\u00ed hodnoty, aby se zajistila stabilita a p\u0159esnost v\u00fdsledk\u016f. Uva\u017eovali jsme tak\u00e9 o mo\u017enosti \u010derpat ze \u010dtvrtletn\u00edch dat, av\u0161ak tato varianta byla odm\u00edtnuta kv\u016fli sez\u00f3nn\u00edm v\u00fdkyv\u016fm v cen\u00e1ch potravin. Pou\u017eit\u00ed ro\u010dn\u00edch pr\u016fm\u011br\u016f t\u00edmto zp\u016fsobem eliminovalo sez\u00f3nn\u00ed vlivy a poskytlo \u010dist\u0161\u00ed data pro anal\u00fdzu dlouhodob\u00fdch trend\u016f.\"\n$newSecondParagraphText = \"Tento p\u0159\u00edstup n\u00e1m umo\u017enil l\u00e9pe porozum\u011bt dynamice v\u00fdvoje mezd a cen potravin v kontextu ekonomick\u00e9ho prost\u0159ed\u00ed \u010cesk\u00e9 republiky.\"\n\n# Locate the target paragraph robustly (by its distinctive original text)\n# instead of assuming a fixed index.\n$targetIndex = 0\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains($marker)) {\n        $targetIndex = $p.Index\n        break\n    }\n}\nif ($targetIndex -eq 0) {\n    throw \"Could not locate the target paragraph (text starting with '$marker').\"\n}\n\n# Replace the whole (multi-run) paragraph text with the revised wording in a\n# single operation - Find/Replace collapses the matched range into one run,\n# mirroring the merge of the four original runs into one.\n$found = $d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newFirstParagraphText, 2)\n\n# Re-fetch the (now single-run) paragraph so we can append new paragraphs\n# right after it.\n$target = $d.Paragraphs($targetIndex)\n\n# Insert the new concluding paragraph right after it.\n$target.Range.InsertParagraphAfter()\n$secondParagraph = $d.Paragraphs($targetIndex + 1)\n$secondParagraph.Range.InsertAfter($newSecondParagraphText)\n\n# Insert a trailing empty paragraph right after that one.\n$secondParagraph.Range.InsertParagraphAfter()\n"}
